# Edit script: Tue, May 12, 2020  6:11:53 PM
#
# 1) Slide 5's table (2nd shape on that slide) switches from the deck's
#    custom "Table_0" style to the built-in table style
#    {B40145A0-20C7-4FF7-B2A5-3C5F7D466B1D}.
# 2) The presentation's theme (carried on the slide master / theme1.xml)
#    switches its colour scheme from the "Integral" / "Red Violet" palette
#    to the stock "Office" palette.

$p = $ppt.ActivePresentation

# -- 1. Table style on slide 5 -------------------------------------------
$tableShape = $p.Slides.Item(5).Shapes.Item(2)
$tableShape.Table.ApplyStyle("{B40145A0-20C7-4FF7-B2A5-3C5F7D466B1D}")

# -- 2. Theme colours -----------------------------------------------------
# ThemeColorScheme.Colors(i) is 1-based in the fixed order:
#   1 dk1, 2 lt1, 3 dk2, 4 lt2, 5-10 accent1-6, 11 hlink, 12 folHlink
# RGB is encoded the usual COM way: R + G*256 + B*65536.
$officeColors = @(
    0,           # dk1      000000
    16777215,    # lt1      FFFFFF
    6968388,     # dk2      44546A
    15132391,    # lt2      E7E6E6
    13998939,    # accent1  5B9BD5
    3243501,     # accent2  ED7D31
    10855845,    # accent3  A5A5A5
    49407,       # accent4  FFC000
    12874308,    # accent5  4472C4
    4697456,     # accent6  70AD47
    12673797,    # hlink    0563C1
    7491477      # folHlink 954F72
)

$themeColors = $p.SlideMaster.Theme.ThemeColorScheme
for ($i = 1; $i -le $officeColors.Length; $i++) {
    $themeColors.Colors($i).RGB = $officeColors[$i - 1]
}
